$wb = $excel.ActiveWorkbook

$values = @(0,0,56796000,-155392000,0,107845000,-230965000,-110107000,0,103783000,147393000,0,-175207000,0,197090000,154862000,174030000,0,252217000,72826000,139623000,0,140988000,82781000,105534000,0,213000000,151000000,163000000,0,245000000,176000000,76000000,0,639000000,245000000,317000000,0,895000000,554000000,688000000,0,1058000000,1295000000,1157000000,0,927000000,358000000,571000000,0,1398000000,976000000,651000000,0,2671000000,2237000000,2444000000,0,601000000,2202000000,499000000,0,10417000000,2524000000,6800000000,0,21869000000,17319000000,18642000000,0,22249000000,28440000000)

# DATA_RAW sheet: new row 24
$ws1 = $wb.Worksheets.Item("DATA_RAW")
$row1 = 24
$ws1.Cells.Item($row1, 1).Value = "FAVÖK"
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws1.Cells.Item($row1, $i + 2).Value = $values[$i]
}

# gelir tablosu (çeyreklik) sheet: new row 11
$ws3 = $wb.Worksheets.Item("gelir tablosu (çeyreklik)")
$row3 = 11
$ws3.Cells.Item($row3, 1).Value = "FAVÖK"
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws3.Cells.Item($row3, $i + 2).Value = $values[$i]
}
